{"js": "// Update the two-digit division answer table: each data row (0, 4, 8, 12, 16)\n// holds 5 \"A\u00f7B=C, D\" problems in its cells; replace each with its new value.\n// Cell (row, col) -> new text, in document order (matches the unified diff).\nconst replacements = [\n  [0, 0, \"39\u00f74=9, 3\"],\n  [0, 1, \"92\u00f77=13, 1\"],\n  [0, 2, \"45\u00f76=7, 3\"],\n  [0, 3, \"42\u00f73=14, 0\"],\n  [0, 4, \"86\u00f73=28, 2\"],\n  [4, 0, \"26\u00f74=6, 2\"],\n  [4, 1, \"19\u00f78=2, 3\"],\n  [4, 2, \"21\u00f79=2, 3\"],\n  [4, 3, \"39\u00f75=7, 4\"],\n  [4, 4, \"40\u00f76=6, 4\"],\n  [8, 0, \"81\u00f76=13, 3\"],\n  [8, 1, \"73\u00f77=10, 3\"],\n  [8, 2, \"30\u00f79=3, 3\"],\n  [8, 3, \"93\u00f76=15, 3\"],\n  [8, 4, \"75\u00f79=8, 3\"],\n  [12, 0, \"49\u00f76=8, 1\"],\n  [12, 1, \"83\u00f72=41, 1\"],\n  [12, 2, \"95\u00f74=23, 3\"],\n  [12, 3, \"78\u00f72=39, 0\"],\n  [12, 4, \"74\u00f73=24, 2\"],\n  [16, 0, \"56\u00f79=6, 2\"],\n  [16, 1, \"85\u00f77=12, 1\"],\n  [16, 2, \"57\u00f76=9, 3\"],\n  [16, 3, \"16\u00f77=2, 2\"],\n  [16, 4, \"99\u00f75=19, 4\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, text] of replacements) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit division answer table: each data row (1, 5, 9, 13, 17 in\n# 1-based Word COM indexing) holds 5 \"A\u00f7B=C, D\" problems in its cells; replace\n# each with its new value. (row, col, newText) in document order, matching the\n# unified diff.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, \"39\u00f74=9, 3\"),\n    @(1, 2, \"92\u00f77=13, 1\"),\n    @(1, 3, \"45\u00f76=7, 3\"),\n    @(1, 4, \"42\u00f73=14, 0\"),\n    @(1, 5, \"86\u00f73=28, 2\"),\n    @(5, 1, \"26\u00f74=6, 2\"),\n    @(5, 2, \"19\u00f78=2, 3\"),\n    @(5, 3, \"21\u00f79=2, 3\"),\n    @(5, 4, \"39\u00f75=7, 4\"),\n    @(5, 5, \"40\u00f76=6, 4\"),\n    @(9, 1, \"81\u00f76=13, 3\"),\n    @(9, 2, \"73\u00f77=10, 3\"),\n    @(9, 3, \"30\u00f79=3, 3\"),\n    @(9, 4, \"93\u00f76=15, 3\"),\n    @(9, 5, \"75\u00f79=8, 3\"),\n    @(13, 1, \"49\u00f76=8, 1\"),\n    @(13, 2, \"83\u00f72=41, 1\"),\n    @(13, 3, \"95\u00f74=23, 3\"),\n    @(13, 4, \"78\u00f72=39, 0\"),\n    @(13, 5, \"74\u00f73=24, 2\"),\n    @(17, 1, \"56\u00f79=6, 2\"),\n    @(17, 2, \"85\u00f77=12, 1\"),\n    @(17, 3, \"57\u00f76=9, 3\"),\n    @(17, 4, \"16\u00f77=2, 2\"),\n    @(17, 5, \"99\u00f75=19, 4\")\n)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $text = $entry[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $text\n}\n"}
